$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 4272.1177  # was 4401.121
$ws.Range("I11").Value = 4272.1177  # was 4401.121
$ws.Range("K11").Value = 4272.1177  # was 4401.121
$ws.Range("M11").Value = -4132.1177  # was -4261.121
$ws.Range("H48").Value = 5716.7144  # was 15005.667
$ws.Range("J48").Value = 4166.6665  # was 15000
$ws.Range("L48").Value = 12499.9995  # was 45000
$ws.Range("N48").Value = -13083.9995  # was -45584
$ws.Range("H56").Value = 5716.7144  # was 15005.667
$ws.Range("J56").Value = 4166.6665  # was 15000
$ws.Range("L56").Value = 12499.9995  # was 45000
$ws.Range("N56").Value = -13567.9995  # was -46068
$ws.Range("H59").Value = 3500  # was 5000
$ws.Range("J59").Value = 3500  # was 5000
$ws.Range("L59").Value = 10500  # was 15000
$ws.Range("N59").Value = -11614  # was -16114
$ws.Range("H98").Value = 10895.223  # was 12186.375
$ws.Range("I98").Value = 10895.223  # was 12186.375
$ws.Range("K98").Value = 10895.223  # was 12186.375
$ws.Range("M98").Value = -9397.223  # was -10688.375
$ws.Range("H122").Value = 10895.223  # was 12186.375
$ws.Range("I122").Value = 10895.223  # was 12186.375
$ws.Range("K122").Value = 32685.669  # was 36559.125
$ws.Range("M122").Value = -30235.669  # was -34109.125

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H30").Value = 2501.75  # was 2123.2
$ws.Range("I30").Value = 1254  # was 1039
$ws.Range("K30").Value = 1254  # was 1039
$ws.Range("M30").Value = -1104  # was -889
$ws.Range("H124").Value = 50635  # was 54037
$ws.Range("J124").Value = 50635  # was 54037
$ws.Range("L124").Value = 50635  # was 54037
$ws.Range("N124").Value = -60455  # was -63857

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H74").Value = 17657  # was 0
$ws.Range("I74").Value = 5000  # was 0
$ws.Range("J74").Value = 30314  # was 0
$ws.Range("K74").Value = 5000  # was 0
$ws.Range("L74").Value = 30314  # was 0
$ws.Range("M74").Value = -4126  # was None
$ws.Range("N74").Value = -32062  # was None
$ws.Range("H77").Value = 17657  # was 0
$ws.Range("I77").Value = 5000  # was 0
$ws.Range("J77").Value = 30314  # was 0
$ws.Range("K77").Value = 15000  # was 0
$ws.Range("L77").Value = 90942  # was 0
$ws.Range("M77").Value = -10632  # was None
$ws.Range("N77").Value = -99678  # was None
$ws.Range("H132").Value = 53696.25  # was 56024.348
$ws.Range("I132").Value = 1818.7142  # was 2005.1666
$ws.Range("J132").Value = 126324.8  # was 114954.37
$ws.Range("K132").Value = 5456.142599999999  # was 6015.4998
$ws.Range("L132").Value = 378974.4  # was 344863.11
$ws.Range("M132").Value = -2926.142599999999  # was -3485.4998
$ws.Range("N132").Value = -384034.4  # was -349923.11

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 77.2  # was 85.22221999999999
$ws.Range("I38").Value = 77.2  # was 85.22221999999999
$ws.Range("K38").Value = 231.6  # was 255.66666
$ws.Range("M38").Value = 115.4  # was 91.33334000000002
$ws.Range("H43").Value = 0  # was 2500
$ws.Range("J43").Value = 0  # was 2500
$ws.Range("L43").Value = 0  # was 7500
$ws.Range("N43").ClearContents()  # was -7728
$ws.Range("H61").Value = 113.833336  # was 110.10526
$ws.Range("J61").Value = 436.33334  # was 338
$ws.Range("L61").Value = 1309.00002  # was 1014
$ws.Range("N61").Value = -1739.00002  # was -1444
$ws.Range("H75").Value = 4862.3335  # was 6632.5
$ws.Range("I75").Value = 1350  # was 1500
$ws.Range("J75").Value = 6618.5  # was 8343.333000000001
$ws.Range("K75").Value = 4050  # was 4500
$ws.Range("L75").Value = 19855.5  # was 25029.999
$ws.Range("M75").Value = -3052  # was -3502
$ws.Range("N75").Value = -21851.5  # was -27025.999
$ws.Range("H78").Value = 4862.3335  # was 6632.5
$ws.Range("I78").Value = 1350  # was 1500
$ws.Range("J78").Value = 6618.5  # was 8343.333000000001
$ws.Range("K78").Value = 12150  # was 13500
$ws.Range("L78").Value = 59566.5  # was 75089.997
$ws.Range("M78").Value = -7158  # was -8508
$ws.Range("N78").Value = -69550.5  # was -85073.997
$ws.Range("H94").Value = 12837.833  # was 12600
$ws.Range("J94").Value = 15006.75  # was 15333.333
$ws.Range("L94").Value = 45020.25  # was 45999.999
$ws.Range("N94").Value = -46372.25  # was -47351.999
$ws.Range("H99").Value = 325  # was 399.33334
$ws.Range("I99").Value = 0  # was 399.33334
$ws.Range("J99").Value = 325  # was 0
$ws.Range("K99").Value = 0  # was 1198.00002
$ws.Range("L99").Value = 975  # was 0
$ws.Range("M99").ClearContents()  # was 1047.99998
$ws.Range("N99").Value = -5467  # was None
$ws.Range("H100").Value = 10985.333  # was 13028
$ws.Range("I100").Value = 10900  # was 0
$ws.Range("J100").Value = 11028  # was 13028
$ws.Range("K100").Value = 32700  # was 0
$ws.Range("L100").Value = 33084  # was 39084
$ws.Range("M100").Value = -31889  # was None
$ws.Range("N100").Value = -34706  # was -40706
$ws.Range("H103").Value = 11285.857  # was 12666.833
$ws.Range("J103").Value = 14989  # was 18985.334
$ws.Range("L103").Value = 44967  # was 56956.00199999999
$ws.Range("N103").Value = -46725  # was -58714.00199999999
$ws.Range("H104").Value = 10544  # was 10894
$ws.Range("I104").Value = 10544  # was 13858.667
$ws.Range("J104").Value = 0  # was 2000
$ws.Range("K104").Value = 31632  # was 41576.001
$ws.Range("L104").Value = 0  # was 6000
$ws.Range("M104").Value = -29011  # was -38955.001
$ws.Range("N104").ClearContents()  # was -11242
$ws.Range("H109").Value = 5897.7144  # was 5464.7144
$ws.Range("I109").Value = 3063.5  # was 3650.6
$ws.Range("J109").Value = 9676.666999999999  # was 10000
$ws.Range("K109").Value = 9190.5  # was 10951.8
$ws.Range("L109").Value = 29030.001  # was 30000
$ws.Range("M109").Value = -8150.5  # was -9911.799999999999
$ws.Range("N109").Value = -31110.001  # was -32080
$ws.Range("H110").Value = 11597.4  # was 11010
$ws.Range("I110").Value = 5975.6665  # was 6500
$ws.Range("K110").Value = 17926.9995  # was 19500
$ws.Range("M110").Value = -13836.9995  # was -15410
$ws.Range("H112").Value = 2866841.2  # was 3343514.8
$ws.Range("J112").Value = 16222.5  # was 19363.334
$ws.Range("L112").Value = 48667.5  # was 58090.00199999999
$ws.Range("N112").Value = -50883.5  # was -60306.00199999999
$ws.Range("H115").Value = 9838  # was 1285.4
$ws.Range("I115").Value = 3042.6667  # was 1285.4
$ws.Range("J115").Value = 20031  # was 0
$ws.Range("K115").Value = 9128.000100000001  # was 3856.2
$ws.Range("L115").Value = 60093  # was 0
$ws.Range("M115").Value = -7953.000100000001  # was -2681.2
$ws.Range("N115").Value = -62443  # was None
$ws.Range("H120").Value = 15806.4  # was 18333.334
$ws.Range("I120").Value = 9999.5  # was 10000
$ws.Range("J120").Value = 19677.666  # was 22500
$ws.Range("K120").Value = 29998.5  # was 30000
$ws.Range("L120").Value = 59032.99800000001  # was 67500
$ws.Range("M120").Value = -25160.5  # was -25162
$ws.Range("N120").Value = -68708.99800000001  # was -77176
$ws.Range("H121").Value = 17545538  # was 15874585
$ws.Range("J121").Value = 22224146  # was 19609604
$ws.Range("L121").Value = 66672438  # was 58828812
$ws.Range("N121").Value = -66675058  # was -58831432
$ws.Range("H137").Value = 2901.6  # was 2892.3635
$ws.Range("I137").Value = 1369.3334  # was 1573.7142
$ws.Range("K137").Value = 4108.0002  # was 4721.142599999999
$ws.Range("M137").Value = 991.9997999999996  # was 378.8574000000008

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 28499.166  # was 24199
$ws.Range("J18").Value = 21748.75  # was 12331.667
$ws.Range("L18").Value = 21748.75  # was 12331.667
$ws.Range("N18").Value = -22334.75  # was -12917.667
$ws.Range("H69").Value = 49800  # was 32256
$ws.Range("J69").Value = 49800  # was 32256
$ws.Range("L69").Value = 49800  # was 32256
$ws.Range("N69").Value = -51298  # was -33754
$ws.Range("H72").Value = 49800  # was 32256
$ws.Range("J72").Value = 49800  # was 32256
$ws.Range("L72").Value = 149400  # was 96768
$ws.Range("N72").Value = -156888  # was -104256
$ws.Range("H94").Value = 0  # was 22250
$ws.Range("J94").Value = 0  # was 22250
$ws.Range("L94").Value = 0  # was 22250
$ws.Range("N94").ClearContents()  # was -23602

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2374.0833  # was 2787.6667
$ws.Range("I16").Value = 1573.375  # was 1655.2858
$ws.Range("J16").Value = 3975.5  # was 6751
$ws.Range("K16").Value = 1573.375  # was 1655.2858
$ws.Range("L16").Value = 3975.5  # was 6751
$ws.Range("M16").Value = -1403.375  # was -1485.2858
$ws.Range("N16").Value = -4315.5  # was -7091
$ws.Range("H22").Value = 2783.725  # was 2848.6924
$ws.Range("I22").Value = 1637.44  # was 1695.25
$ws.Range("K22").Value = 1637.44  # was 1695.25
$ws.Range("M22").Value = -1342.44  # was -1400.25
$ws.Range("H27").Value = 2783.725  # was 2848.6924
$ws.Range("I27").Value = 1637.44  # was 1695.25
$ws.Range("K27").Value = 1637.44  # was 1695.25
$ws.Range("M27").Value = -1530.44  # was -1588.25
$ws.Range("H122").Value = 5311.8237  # was 5573.467
$ws.Range("I122").Value = 4409.9  # was 4675
$ws.Range("K122").Value = 13229.7  # was 14025
$ws.Range("M122").Value = -10779.7  # was -11575
$ws.Range("H132").Value = 9407.066000000001  # was 9066.375
$ws.Range("I132").Value = 8937.888999999999  # was 8449.700000000001
$ws.Range("J132").Value = 10110.833  # was 10094.167
$ws.Range("K132").Value = 26813.667  # was 25349.1
$ws.Range("L132").Value = 30332.499  # was 30282.501
$ws.Range("M132").Value = -24283.667  # was -22819.1
$ws.Range("N132").Value = -35392.499  # was -35342.501

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 13916.667  # was 14750
$ws.Range("I5").Value = 500  # was 5000
$ws.Range("J5").Value = 16600  # was 16700
$ws.Range("K5").Value = 500  # was 5000
$ws.Range("L5").Value = 16600  # was 16700
$ws.Range("M5").Value = -388  # was -4888
$ws.Range("N5").Value = -16824  # was -16924
$ws.Range("H15").Value = 50000000  # was 16685002
$ws.Range("I15").Value = 0  # was 20000
$ws.Range("J15").Value = 50000000  # was 25017504
$ws.Range("K15").Value = 0  # was 20000
$ws.Range("L15").Value = 50000000  # was 25017504
$ws.Range("M15").ClearContents()  # was -19712
$ws.Range("N15").Value = -50000576  # was -25018080
$ws.Range("H81").Value = 8388.666999999999  # was 8937.25
$ws.Range("I81").Value = 3249.75  # was 2999.6667
$ws.Range("K81").Value = 6499.5  # was 5999.3334
$ws.Range("M81").Value = -5438.5  # was -4938.3334
$ws.Range("H84").Value = 8388.666999999999  # was 8937.25
$ws.Range("I84").Value = 3249.75  # was 2999.6667
$ws.Range("K84").Value = 32497.5  # was 29996.667
$ws.Range("M84").Value = -27193.5  # was -24692.667
$ws.Range("H92").Value = 0  # was 20000
$ws.Range("J92").Value = 0  # was 20000
$ws.Range("L92").Value = 0  # was 20000
$ws.Range("N92").ClearContents()  # was -24992
$ws.Range("H100").Value = 5026.3076  # was 26000
$ws.Range("I100").Value = 5278.5  # was 50000
$ws.Range("K100").Value = 10557  # was 100000
$ws.Range("M100").Value = -10016  # was -99459
$ws.Range("H113").Value = 2056.5557  # was 2497.7144
$ws.Range("I113").Value = 1001.6667  # was 1330
$ws.Range("J113").Value = 4166.3335  # was 3373.5
$ws.Range("K113").Value = 3005.0001  # was 3990
$ws.Range("L113").Value = 12499.0005  # was 10120.5
$ws.Range("M113").Value = -835.0001000000002  # was -1820
$ws.Range("N113").Value = -16839.0005  # was -14460.5
$ws.Range("H132").Value = 5097.4053  # was 5371.543
$ws.Range("I132").Value = 2414.6956  # was 2616.0952
$ws.Range("K132").Value = 7244.0868  # was 7848.285600000001
$ws.Range("M132").Value = -4714.0868  # was -5318.285600000001
